$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update status column (C) from "In Progress" to "Complete" for several rows
$completeRows = @(7, 12, 14, 15, 16, 17, 18, 19, 20)
foreach ($r in $completeRows) {
    $ws.Cells.Item($r, 3).Value = "Complete"
}

# 2. Add new "Sprint 3" section starting at row 21
$ws.Cells.Item(21, 1).Value = "Sprint 3"
$ws.Cells.Item(21, 1).HorizontalAlignment = -4108

# 3. Add the three new Sprint 3 tasks (rows 22-24)
$ws.Cells.Item(22, 1).Value = "Task 1"
$ws.Cells.Item(22, 1).HorizontalAlignment = -4108
$ws.Cells.Item(22, 2).Value = "Add menu with game pieces (Free Play Level)"
$ws.Cells.Item(22, 3).Value = "In Progress"

$ws.Cells.Item(23, 1).Value = "Task 2"
$ws.Cells.Item(23, 1).HorizontalAlignment = -4108
$ws.Cells.Item(23, 2).Value = "Add functionality with game pieces on the cell  (Free Play Level)"
$ws.Cells.Item(23, 3).Value = "In Progress"

$ws.Cells.Item(24, 1).Value = "Task 3"
$ws.Cells.Item(24, 1).HorizontalAlignment = -4108
$ws.Cells.Item(24, 2).Value = "Allow game pieces to interact when the play button is clicked  (Free Play Level)"
$ws.Cells.Item(24, 3).Value = "In Progress"

# 4. Add trailing blank/centered rows 25-27
$ws.Cells.Item(25, 1).HorizontalAlignment = -4108
$ws.Cells.Item(26, 1).HorizontalAlignment = -4108
$ws.Cells.Item(27, 1).HorizontalAlignment = -4108

# 5. Update the active selection to reflect the new working cell
$ws.Range("G20").Select()
